# Update "Provincias Spain" workbook: refresh the timestamp banner and the
# per-province COVID figures (Casos totales / Casos activos / Recuperados / Muertes)
# for the rows whose numbers moved in this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header banner timestamp update (row 1 / cell A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 12:29"

# Row 4 - Madrid
$ws.Range("B4").Value = 21520
$ws.Range("C4").Value = 6326
$ws.Range("D4").Value = 12437
$ws.Range("E4").Value = 2757

# Row 7 - Bizkaia/Vizcaya
$ws.Range("B7").Value = 2263
$ws.Range("C7").Value = 68
$ws.Range("D7").Value = 2091
$ws.Range("E7").Value = 104

# Row 8 - Valencia/Valencia
$ws.Range("B8").Value = 2186
$ws.Range("C8").Value = 814
$ws.Range("D8").Value = 1793
$ws.Range("E8").Value = 79

# Row 9 - Navarra
$ws.Range("B9").Value = 1829
$ws.Range("C9").Value = 98
$ws.Range("D9").Value = 1661
$ws.Range("E9").Value = 70

# Row 11 - Ciudad Real
$ws.Range("B11").Value = 1436
$ws.Range("C11").Value = 364
$ws.Range("D11").Value = 1007
$ws.Range("E11").Value = 65

# Row 12 - La Rioja
$ws.Range("B12").Value = 1422
$ws.Range("C12").Value = 153
$ws.Range("D12").Value = 1273
$ws.Range("E12").Value = 89

# Row 13 - Toledo
$ws.Range("B13").Value = 1322
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 1193
$ws.Range("E13").Value = 109

# Row 14 - Alacant/Alicante
$ws.Range("B14").Value = 1112
$ws.Range("C14").Value = 153
$ws.Range("D14").Value = 972
$ws.Range("E14").Value = 90

# Row 20 - Aragon
$ws.Range("B20").Value = 937
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 894
$ws.Range("E20").Value = 22

# Row 21 - Caceres
$ws.Range("B21").Value = 907
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 838
$ws.Range("E21").Value = 40

# Row 22 - Cantabria
$ws.Range("B22").Value = 841
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 777
$ws.Range("E22").Value = 60

# Row 24 - Murcia
$ws.Range("B24").Value = 802
$ws.Range("D24").Value = 773

# Row 34 - Guadalajara
$ws.Range("B34").Value = 449
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 424
$ws.Range("E34").Value = 21

# Row 35 - Jaen
$ws.Range("B35").Value = 440
$ws.Range("C35").Value = 153
$ws.Range("D35").Value = 362
$ws.Range("E35").Value = 75

# Row 36 - Castello/Castellon
$ws.Range("B36").Value = 414
$ws.Range("C36").Value = 7
$ws.Range("D36").Value = 392
$ws.Range("E36").Value = 15

# Row 54 - Melilla
$ws.Range("B54").Value = 45
$ws.Range("D54").Value = 44
